$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Smith, John"
$ws.Range("D2").Value = 37092
$ws.Range("D2").NumberFormat = "#,##0"
$ws.Range("Q2").Value = 575
$ws.Range("Q2").NumberFormat = "0.00"
$ws.Range("R2").Value = "LPB"
$ws.Range("T2").Value = 75
$ws.Range("T2").NumberFormat = "0"

# ---------------------------------------------------------------------------
# Row 3
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "Nguyen, Kim-Vincent"
$ws.Range("C3").Value = "1ABC124"
$ws.Range("D3").Value = 39769
$ws.Range("D3").NumberFormat = "#,##0"
$ws.Range("G3").Value = "Rental"
$ws.Range("H3").Value = "Diesel"
$ws.Range("I3").Value = "BMW"
$ws.Range("J3").Value = "SERIE 1"
$ws.Range("K3").Value = "M135i 306 "
$ws.Range("Q3").Value = 575.01
$ws.Range("R3").Value = "LPB"
$ws.Range("T3").Value = 75.010000000000005

# ---------------------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "Nguyen, Khanh-Michel"
$ws.Range("C4").Value = "1ABC125"
$ws.Range("D4").Value = 41254
$ws.Range("D4").NumberFormat = "#,##0"
$ws.Range("H4").Value = "Diesel"
$ws.Range("I4").Value = "AUDI"
$ws.Range("J4").Value = "A3"
$ws.Range("K4").Value = "2.5 L TFSI 400"
$ws.Range("Q4").Value = 575.02
$ws.Range("R4").Value = "LPB"
$ws.Range("T4").Value = 75.02

# ---------------------------------------------------------------------------
# Row 5
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "Wick, John"
$ws.Range("C5").Value = "1ABC126"
$ws.Range("D5").Value = 47749
$ws.Range("D5").NumberFormat = "#,##0"
$ws.Range("I5").Value = "MERCEDES"
$ws.Range("J5").Value = "A-CLASS"
$ws.Range("K5").Value = "250 4MATIC BlueEfficiency"
$ws.Range("Q5").Value = 575.03
$ws.Range("R5").Value = "LPB"
$ws.Range("T5").Value = 75.03

# ---------------------------------------------------------------------------
# Row 6
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = "Belcaid, Younes"
$ws.Range("C6").Value = "1ABC127"
$ws.Range("D6").Value = 42835
$ws.Range("D6").NumberFormat = "#,##0"
$ws.Range("G6").Value = "Rental"
$ws.Range("I6").Value = "MERCEDES"
$ws.Range("J6").Value = "A-CLASS"
$ws.Range("K6").Value = "A 160 CDI"
$ws.Range("Q6").Value = 575.04
$ws.Range("R6").Value = "LPB"
$ws.Range("T6").Value = 75.040000000000006

# ---------------------------------------------------------------------------
# Column width (Tax value column)
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 11.14

# ---------------------------------------------------------------------------
# View settings: zoom + selection on the active window
# ---------------------------------------------------------------------------
$win = $ws.Application.ActiveWindow
$win.Zoom = 130
$ws.Range("O2").Select()
